$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main Site")

# New header cell for the desktopOnly column
$ws.Range("E1").Value = "desktopOnly"
$ws.Range("E1").Font.Bold = $true

# Set column width to match diff (stored worksheet width 15.86328125 chars;
# closest value reachable via the ColumnWidth COM property given this
# engine's pixel-quantized width storage)
$ws.Columns.Item(5).ColumnWidth = 15

# Boolean value for row 3
$ws.Range("E3").Value = $true

# Move selection to E3, matching the diff
$ws.Range("E3").Select()
